# Apply the "Add units to forms and fix logic for grid availability visibility" edit.

$wb = $excel.ActiveWorkbook

$survey  = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# ---------------------------------------------------------------------------
# 1. survey sheet: insert two new columns (J, K) for display.hint.text /
#    display.hint.text.es, pushing the old choice_filter/required columns
#    from J/K to L/M.
# ---------------------------------------------------------------------------
$survey.Range("J1:K1").EntireColumn.Insert()

# New header row cells.
$survey.Range("J1").Value = "display.hint.text"
$survey.Range("K1").Value = "display.hint.text.es"

# Unit hints for "Grid Availability" (hours per day) ...
$survey.Range("J18").Value = "(hours per day)"
$survey.Range("K18").Value = "(horas por día)"

# ... and for "Distance to the Closest Supply Point" (kilometers).
$survey.Range("J23").Value = "(kilometers)"
$survey.Range("K23").Value = "(kilómetros)"

# ---------------------------------------------------------------------------
# 2. Fix the grid-availability visibility condition so it also triggers for
#    the split "both_grid_and_generator" / "both_grid_and_solar" choices
#    (previously it only checked for the old combined "both" choice).
# ---------------------------------------------------------------------------
$survey.Range("B17").Value = "(selected(data('electricity_source'), 'grid') || selected(data('electricity_source'), 'both_grid_and_generator') || selected(data('electricity_source'), 'both_grid_and_solar'))"

# ---------------------------------------------------------------------------
# 3. View-state tweaks: survey becomes the active/selected sheet, choices
#    loses it; update each sheet's current selection.
# ---------------------------------------------------------------------------
$survey.Activate()
$survey.Range("K30").Select()

$choices.Range("B30").Select()

$survey.Activate()
